{"js": "// Replace the date line and each two-digit multiplication problem's text\n// with its new value, preserving all existing run/paragraph formatting by\n// using search + replace (InsertLocation.Replace) instead of rebuilding\n// the runs from scratch.\nconst replacements = [\n  [\"2023-08-04 Friday\", \"2023-08-05 Saturday\"],\n  [\"78\u00d743=\", \"67\u00d761=\"],\n  [\"93\u00d763=\", \"65\u00d798=\"],\n  [\"59\u00d717=\", \"18\u00d713=\"],\n  [\"27\u00d741=\", \"53\u00d724=\"],\n  [\"99\u00d789=\", \"22\u00d791=\"],\n  [\"85\u00d785=\", \"94\u00d714=\"],\n  [\"88\u00d747=\", \"88\u00d787=\"],\n  [\"53\u00d732=\", \"36\u00d727=\"],\n  [\"71\u00d745=\", \"96\u00d757=\"],\n  [\"71\u00d775=\", \"84\u00d711=\"],\n  [\"59\u00d790=\", \"67\u00d722=\"],\n  [\"14\u00d796=\", \"21\u00d733=\"],\n  [\"64\u00d731=\", \"75\u00d749=\"],\n  [\"31\u00d754=\", \"62\u00d731=\"],\n  [\"80\u00d742=\", \"83\u00d753=\"],\n  [\"90\u00d729=\", \"60\u00d738=\"],\n  [\"61\u00d720=\", \"47\u00d728=\"],\n  [\"55\u00d750=\", \"56\u00d763=\"],\n  [\"29\u00d714=\", \"20\u00d752=\"],\n  [\"14\u00d713=\", \"19\u00d712=\"],\n  [\"50\u00d765=\", \"41\u00d796=\"],\n  [\"71\u00d766=\", \"64\u00d752=\"],\n  [\"58\u00d728=\", \"97\u00d733=\"],\n  [\"57\u00d719=\", \"62\u00d717=\"],\n  [\"53\u00d771=\", \"93\u00d746=\"],\n];\n\nfor (const [from, to] of replacements) {\n  const results = context.document.body.search(from, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(to, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and each two-digit multiplication problem's text\n# with its new value, preserving all existing run/paragraph formatting by\n# using Find/Replace (wdReplaceOne) instead of rebuilding the runs.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-08-04 Friday\", \"2023-08-05 Saturday\"),\n    @(\"78\u00d743=\", \"67\u00d761=\"),\n    @(\"93\u00d763=\", \"65\u00d798=\"),\n    @(\"59\u00d717=\", \"18\u00d713=\"),\n    @(\"27\u00d741=\", \"53\u00d724=\"),\n    @(\"99\u00d789=\", \"22\u00d791=\"),\n    @(\"85\u00d785=\", \"94\u00d714=\"),\n    @(\"88\u00d747=\", \"88\u00d787=\"),\n    @(\"53\u00d732=\", \"36\u00d727=\"),\n    @(\"71\u00d745=\", \"96\u00d757=\"),\n    @(\"71\u00d775=\", \"84\u00d711=\"),\n    @(\"59\u00d790=\", \"67\u00d722=\"),\n    @(\"14\u00d796=\", \"21\u00d733=\"),\n    @(\"64\u00d731=\", \"75\u00d749=\"),\n    @(\"31\u00d754=\", \"62\u00d731=\"),\n    @(\"80\u00d742=\", \"83\u00d753=\"),\n    @(\"90\u00d729=\", \"60\u00d738=\"),\n    @(\"61\u00d720=\", \"47\u00d728=\"),\n    @(\"55\u00d750=\", \"56\u00d763=\"),\n    @(\"29\u00d714=\", \"20\u00d752=\"),\n    @(\"14\u00d713=\", \"19\u00d712=\"),\n    @(\"50\u00d765=\", \"41\u00d796=\"),\n    @(\"71\u00d766=\", \"64\u00d752=\"),\n    @(\"58\u00d728=\", \"97\u00d733=\"),\n    @(\"57\u00d719=\", \"62\u00d717=\"),\n    @(\"53\u00d771=\", \"93\u00d746=\")\n)\n\nforeach ($pair in $replacements) {\n    $from = $pair[0]\n    $to = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $from\n    $find.Replacement.Text = $to\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
